# Daily attendance processing - 2025-11-05 04:25:47
# Re-order the "Recorded By" (column G) value lists so that the first
# literal "System" token (exact case) is moved to the end of the
# comma-separated list, leaving everything else (including any
# lower-case "system" entries) in its original relative order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) {
    $lastRow = 1
}

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }

    $parts = $text -split ", "
    $newParts = @()
    $systemToken = $null

    foreach ($p in $parts) {
        if ($p.Equals("System") -and ($systemToken -eq $null)) {
            $systemToken = $p
        } else {
            $newParts += $p
        }
    }

    if ($systemToken -ne $null) {
        $newParts += $systemToken
        $result = $newParts -join ", "
        if (-not $result.Equals($text)) {
            $cell.Value = $result
        }
    }
}
